# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells must remain stored as text, matching the original
# inline-string cell type, so force text number format before assignment and
# restore the default style afterward to avoid leaving a stray numeric format.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D10", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D23", "D27", "D29", "D30", "D33", "D34", "D36", "D38", "D39", "D40", "D42", "D43", "D48", "D50", "D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.060.71"
$ws.Range("D3").Value = "2.372.01"
$ws.Range("D5").Value = "303.25"
$ws.Range("D6").Value = "95.52"
$ws.Range("D7").Value = "0.503"
$ws.Range("D8").Value = "1.00"
$ws.Range("D10").Value = "34.33"
$ws.Range("D13").Value = "18.35"
$ws.Range("D14").Value = "6.74"
$ws.Range("D15").Value = "2.736.06"
$ws.Range("D16").Value = "2.368.50"
$ws.Range("D17").Value = "0.797"
$ws.Range("D18").Value = "43.055.37"
$ws.Range("D19").Value = "11.99"
$ws.Range("D20").Value = "6.28"
$ws.Range("D23").Value = "235.60"
$ws.Range("D27").Value = "24.46"
$ws.Range("D29").Value = "9.32"
$ws.Range("D30").Value = "32.43"
$ws.Range("D33").Value = "17.63"
$ws.Range("D34").Value = "0.0727"
$ws.Range("D36").Value = "129.20"
$ws.Range("D38").Value = "2.86"
$ws.Range("D39").Value = "4.34"
$ws.Range("D40").Value = "2.26"
$ws.Range("D42").Value = "20.94"
$ws.Range("D43").Value = "1.928.52"
$ws.Range("D48").Value = "2.595.70"
$ws.Range("D50").Value = "71.42"
$ws.Range("D51").Value = "1.13"

foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).Style = "Normal"
}

# Remaining text columns (Coin name, Link, Volume) are safe to assign directly.
$ws.Range("E2").Value = "  +0.87%  "
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +1.98%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("E7").Value = "  +0.69%  "
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  -2.24%  "
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("E11").Value = "  +3.10%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("E13").Value = "  -1.14%  "
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  +1.39%  "
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  +2.12%  "
$ws.Range("E21").Value = "  +0.35%  "
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("E24").Value = "  -1.04%  "
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("E28").Value = "  +2.84%  "
$ws.Range("E29").Value = "  +2.65%  "
$ws.Range("E30").Value = "  +4.17%  "
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  -4.07%  "
$ws.Range("E35").Value = "  +7.55%  "
$ws.Range("E36").Value = "  -2.45%  "
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("E38").Value = "  +5.39%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("E45").Value = "  +3.30%  "
$ws.Range("E46").Value = "  -8.71%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("E51").Value = "  +1.23%  "

